# Append two new log rows (156, 157) to the feed_logs sheet, matching
# the run_id/rss_url_id/date/response/item_count columns already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 156
$ws.Cells.Item(156, 1).Value = 155
$ws.Cells.Item(156, 2).Value = 1
$ws.Cells.Item(156, 3).Value = "2024-06-18 06:17:42"
$ws.Cells.Item(156, 4).Value = 200
$ws.Cells.Item(156, 5).Value = 10

# Row 157
$ws.Cells.Item(157, 1).Value = 156
$ws.Cells.Item(157, 2).Value = 2
$ws.Cells.Item(157, 3).Value = "2024-06-18 06:17:43"
$ws.Cells.Item(157, 4).Value = 200
$ws.Cells.Item(157, 5).Value = 0
